$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 11 values
$ws.Range("A11").Value = "DataBase Task (Waleed)"
$ws.Range("G11").Value = "Done"

# Update selection
$ws.Range("C16").Select()
